$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '61.366.83'
Set-TextCell 'E2' '  +7.55%  '

Set-TextCell 'D3' '3.323.10'
Set-TextCell 'E3' '  +2.39%  '

Set-TextCell 'D4' '0.999'
Set-TextCell 'E4' '  -0.22%  '

Set-TextCell 'D5' '410.80'
Set-TextCell 'E5' '  +3.87%  '

Set-TextCell 'D6' '112.74'
Set-TextCell 'E6' '  +4.46%  '

Set-TextCell 'D7' '3.315.53'
Set-TextCell 'E7' '  +2.25%  '

Set-TextCell 'D8' '0.567'
Set-TextCell 'E8' '  -4.10%  '

Set-TextCell 'E9' '  -0.01%  '

Set-TextCell 'D10' '0.622'
Set-TextCell 'E10' '  -0.60%  '

Set-TextCell 'D11' '0.117'
Set-TextCell 'E11' '  +19.28%  '

Set-TextCell 'D12' '38.83'
Set-TextCell 'E12' '  -0.86%  '

Set-TextCell 'E13' '  -0.01%  '

Set-TextCell 'D14' '3.843.03'
Set-TextCell 'E14' '  +2.21%  '

Set-TextCell 'D15' '8.17'
Set-TextCell 'E15' '  -0.32%  '

Set-TextCell 'D16' '19.04'
Set-TextCell 'E16' '  -0.28%  '

Set-TextCell 'D17' '3.320.92'
Set-TextCell 'E17' '  +2.38%  '

Set-TextCell 'D18' '61.017.92'
Set-TextCell 'E18' '  +7.30%  '

Set-TextCell 'D19' '0.987'
Set-TextCell 'E19' '  -3.88%  '

Set-TextCell 'D20' '10.59'
Set-TextCell 'E20' '  -1.91%  '

Set-TextCell 'D21' '0.0000116'
Set-TextCell 'E21' '  +2.07%  '

Set-TextCell 'D22' '3.23'
Set-TextCell 'E22' '  -3.64%  '

Set-TextCell 'D23' '12.36'
Set-TextCell 'E23' '  -4.69%  '

Set-TextCell 'D24' '295.00'
Set-TextCell 'E24' '  -0.16%  '

Set-TextCell 'D25' '73.26'
Set-TextCell 'E25' '  -1.41%  '

Set-TextCell 'D26' '3.08'
Set-TextCell 'E26' '  -2.66%  '

Set-TextCell 'D27' '28.93'
Set-TextCell 'E27' '  +3.98%  '

Set-TextCell 'D28' '4.53'
Set-TextCell 'E28' '  +4.23%  '

Set-TextCell 'D29' '0.173'
Set-TextCell 'E29' '  +2.68%  '

Set-TextCell 'D30' '7.36'
Set-TextCell 'E30' '  +1.17%  '

Set-TextCell 'D31' '7.42'
Set-TextCell 'E31' '  -2.09%  '

Set-TextCell 'E32' '  -0.05%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D33' '0.109'
Set-TextCell 'E33' '  +0.82%  '

$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D34' '11.11'
Set-TextCell 'E34' '  -2.53%  '

Set-TextCell 'D35' '2.47'
Set-TextCell 'E35' '  +15.97%  '

Set-TextCell 'D36' '39.80'
Set-TextCell 'E36' '  +1.82%  '

Set-TextCell 'D37' '0.0480'
Set-TextCell 'E37' '  -0.32%  '

Set-TextCell 'D38' '52.56'
Set-TextCell 'E38' '  +1.82%  '

Set-TextCell 'E39' '  -0.01%  '

Set-TextCell 'D40' '3.03'
Set-TextCell 'E40' '  +4.18%  '

Set-TextCell 'E41' '  -6.03%  '

Set-TextCell 'D42' '135.15'
Set-TextCell 'E42' '  -0.04%  '

Set-TextCell 'D43' '0.120'
Set-TextCell 'E43' '  -2.24%  '

Set-TextCell 'D44' '1.88'
Set-TextCell 'E44' '  -0.28%  '

Set-TextCell 'D45' '0.284'
Set-TextCell 'E45' '  +0.95%  '

Set-TextCell 'D46' '16.23'
Set-TextCell 'E46' '  -4.89%  '

Set-TextCell 'D47' '3.76'
Set-TextCell 'E47' '  -4.44%  '

Set-TextCell 'D48' '2.19'
Set-TextCell 'E48' '  +2.51%  '

Set-TextCell 'D49' '20.84'
Set-TextCell 'E49' '  -6.28%  '

Set-TextCell 'D50' '2.112.41'
Set-TextCell 'E50' '  -2.20%  '

Set-TextCell 'D51' '3.646.54'
Set-TextCell 'E51' '  +2.24%  '
